$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2199
$ws.Range("E2").Value = 31
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = 31
$ws.Range("H2").Value = 32
$ws.Range("I2").Value = 53
$ws.Range("J2").Value = -20
$ws.Range("K2").Value = 3466
$ws.Range("L2").Value = 1953
$ws.Range("M2").Value = 1514
$ws.Range("N2").Value = 1471
$ws.Range("O2").Value = 43
$ws.Range("P2").Value = 295
$ws.Range("Q2").Value = -46
$ws.Range("R2").Value = 230
$ws.Range("S2").Value = -250
$ws.Range("T2").Value = 16
$ws.Range("U2").Value = -61
$ws.Range("V2").Value = 1406
$ws.Range("W2").Value = 1.42
$ws.Range("X2").Value = 1.47
$ws.Range("Y2").Value = 3.54
$ws.Range("Z2").Value = 0.9
$ws.Range("AA2").Value = 128.99
$ws.Range("AB2").Value = 301.79
$ws.Range("AC2").Value = 78
$ws.Range("AD2").Value = 13.13
$ws.Range("AE2").Value = 2180
$ws.Range("AF2").Value = 0.47
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").ClearContents()
$ws.Range("AJ2").Value = 67450655
$ws.Range("D3").Value = 2091
$ws.Range("E3").Value = 49
$ws.Range("F3").Value = 49
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = -4
$ws.Range("K3").Value = 3405
$ws.Range("L3").Value = 1895
$ws.Range("M3").Value = 1510
$ws.Range("N3").Value = 1467
$ws.Range("O3").Value = 42
$ws.Range("P3").Value = 295
$ws.Range("Q3").Value = 67
$ws.Range("R3").Value = -91
$ws.Range("S3").Value = 9
$ws.Range("T3").Value = 121
$ws.Range("U3").Value = -54
$ws.Range("V3").Value = 1362
$ws.Range("W3").Value = 2.32
$ws.Range("X3").Value = 0.25
$ws.Range("Y3").Value = 0.6
$ws.Range("Z3").Value = 0.15
$ws.Range("AA3").Value = 125.49
$ws.Range("AB3").Value = 302.33
$ws.Range("AC3").Value = 13
$ws.Range("AD3").Value = 75.09
$ws.Range("AE3").Value = 2175
$ws.Range("AF3").Value = 0.45
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").ClearContents()
$ws.Range("AJ3").Value = 67450655
$ws.Range("D4").Value = 2277
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = 40
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = -2
$ws.Range("K4").Value = 3831
$ws.Range("L4").Value = 1809
$ws.Range("M4").Value = 2022
$ws.Range("N4").Value = 1972
$ws.Range("O4").Value = 50
$ws.Range("P4").Value = 445
$ws.Range("Q4").Value = 310
$ws.Range("R4").Value = -291
$ws.Range("S4").Value = 71
$ws.Range("T4").Value = 149
$ws.Range("U4").Value = 161
$ws.Range("V4").Value = 1235
$ws.Range("W4").Value = 1.76
$ws.Range("X4").Value = 0.03
$ws.Range("Y4").Value = 0.17
$ws.Range("Z4").Value = 0.02
$ws.Range("AA4").Value = 89.47
$ws.Range("AB4").Value = 218.43
$ws.Range("AC4").Value = 4
$ws.Range("AD4").Value = 264.93
$ws.Range("AE4").Value = 2215
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").ClearContents()
$ws.Range("AJ4").Value = 89047384
$ws.Range("D5").Value = 2263
$ws.Range("E5").Value = 66
$ws.Range("F5").Value = 66
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = 14
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 3598
$ws.Range("L5").Value = 1372
$ws.Range("M5").Value = 2226
$ws.Range("N5").Value = 2174
$ws.Range("O5").Value = 53
$ws.Range("P5").Value = 554
$ws.Range("Q5").Value = 124
$ws.Range("R5").Value = 28
$ws.Range("S5").Value = -220
$ws.Range("T5").Value = 53
$ws.Range("U5").Value = 71
$ws.Range("V5").Value = 811
$ws.Range("W5").Value = 2.91
$ws.Range("X5").Value = 0.6
$ws.Range("Y5").Value = 0.45
$ws.Range("Z5").Value = 0.37
$ws.Range("AA5").Value = 61.63
$ws.Range("AB5").Value = 188.94
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 130.29
$ws.Range("AE5").Value = 1961
$ws.Range("AF5").Value = 0.66
$ws.Range("AG5").Value = 12
$ws.Range("AH5").Value = 0.93
$ws.Range("AI5").ClearContents()
$ws.Range("AJ5").Value = 110855877
$ws.Range("D6").Value = 2174
$ws.Range("E6").Value = 33
$ws.Range("F6").Value = 33
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = 18
$ws.Range("I6").Value = 15
$ws.Range("K6").Value = 3512
$ws.Range("L6").Value = 1258
$ws.Range("M6").Value = 2254
$ws.Range("N6").Value = 2196
$ws.Range("P6").Value = 562
$ws.Range("Q6").Value = 280
$ws.Range("R6").Value = -97
$ws.Range("S6").Value = -96
$ws.Range("T6").Value = 49
$ws.Range("U6").Value = 231
$ws.Range("V6").Value = 716
$ws.Range("W6").Value = 1.54
$ws.Range("X6").Value = 0.8100000000000001
$ws.Range("Y6").Value = 0.67
$ws.Range("Z6").Value = 0.49
$ws.Range("AA6").Value = 55.8
$ws.Range("AB6").Value = 187.33
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 79.41
$ws.Range("AE6").Value = 1953
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").ClearContents()
$ws.Range("AJ6").Value = 112451621
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
